$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old rows 28:29 that are being removed (fewer data rows now)
$ws.Range("A28:F29").Clear()

# Apply header style (bold, bordered, centered) to new column F1 by copying from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Update headers
$ws.Range("A1").Value = "ID"
$ws.Range("D1").Value = "Quantidade"
$ws.Range("E1").Value = "Lote_id"
$ws.Range("F1").Value = "Max"

# Update data rows (Lote ID -> Material+Qualidade grouping key; add Lote_id sequence and Max column)
$ws.Range("A2").Value = "495223 Q2"
$ws.Range("B2").Value = 495223
$ws.Range("C2").Value = "Q2"
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 25
$ws.Range("A3").Value = "495223 Q2"
$ws.Range("B3").Value = 495223
$ws.Range("C3").Value = "Q2"
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 25
$ws.Range("A4").Value = "495223 Q2"
$ws.Range("B4").Value = 495223
$ws.Range("C4").Value = "Q2"
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 25
$ws.Range("A5").Value = "495223 Q2"
$ws.Range("B5").Value = 495223
$ws.Range("C5").Value = "Q2"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 25
$ws.Range("A6").Value = "495223 Q3"
$ws.Range("B6").Value = 495223
$ws.Range("C6").Value = "Q3"
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 25
$ws.Range("A7").Value = "495223 Q3"
$ws.Range("B7").Value = 495223
$ws.Range("C7").Value = "Q3"
$ws.Range("D7").Value = 25
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 25
$ws.Range("A8").Value = "495223 Q3"
$ws.Range("B8").Value = 495223
$ws.Range("C8").Value = "Q3"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 25
$ws.Range("A9").Value = "495224 Q1"
$ws.Range("B9").Value = 495224
$ws.Range("C9").Value = "Q1"
$ws.Range("D9").Value = 25
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 28
$ws.Range("A10").Value = "495224 Q3"
$ws.Range("B10").Value = 495224
$ws.Range("C10").Value = "Q3"
$ws.Range("D10").Value = 28
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 28
$ws.Range("A11").Value = "495224 Q3"
$ws.Range("B11").Value = 495224
$ws.Range("C11").Value = "Q3"
$ws.Range("D11").Value = 16
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 28
$ws.Range("A12").Value = "495225 Q2"
$ws.Range("B12").Value = 495225
$ws.Range("C12").Value = "Q2"
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 30
$ws.Range("A13").Value = "495225 Q2"
$ws.Range("B13").Value = 495225
$ws.Range("C13").Value = "Q2"
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 30
$ws.Range("A14").Value = "495225 Q2"
$ws.Range("B14").Value = 495225
$ws.Range("C14").Value = "Q2"
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 30
$ws.Range("A15").Value = "495225 Q2"
$ws.Range("B15").Value = 495225
$ws.Range("C15").Value = "Q2"
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 30
$ws.Range("A16").Value = "495225 Q3"
$ws.Range("B16").Value = 495225
$ws.Range("C16").Value = "Q3"
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 30
$ws.Range("A17").Value = "495225 Q3"
$ws.Range("B17").Value = 495225
$ws.Range("C17").Value = "Q3"
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 30
$ws.Range("A18").Value = "495226 Q1"
$ws.Range("B18").Value = 495226
$ws.Range("C18").Value = "Q1"
$ws.Range("D18").Value = 27
$ws.Range("E18").Value = 17
$ws.Range("F18").Value = 27
$ws.Range("A19").Value = "495226 Q1"
$ws.Range("B19").Value = 495226
$ws.Range("C19").Value = "Q1"
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = 18
$ws.Range("F19").Value = 27
$ws.Range("A20").Value = "495226 Q1"
$ws.Range("B20").Value = 495226
$ws.Range("C20").Value = "Q1"
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 19
$ws.Range("F20").Value = 27
$ws.Range("A21").Value = "495226 Q1"
$ws.Range("B21").Value = 495226
$ws.Range("C21").Value = "Q1"
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 27
$ws.Range("A22").Value = "495226 Q2"
$ws.Range("B22").Value = 495226
$ws.Range("C22").Value = "Q2"
$ws.Range("D22").Value = 27
$ws.Range("E22").Value = 21
$ws.Range("F22").Value = 27
$ws.Range("A23").Value = "495226 Q3"
$ws.Range("B23").Value = 495226
$ws.Range("C23").Value = "Q3"
$ws.Range("D23").Value = 27
$ws.Range("E23").Value = 22
$ws.Range("F23").Value = 27
$ws.Range("A24").Value = "495226 Q3"
$ws.Range("B24").Value = 495226
$ws.Range("C24").Value = "Q3"
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 23
$ws.Range("F24").Value = 27
$ws.Range("A25").Value = "495226 Q3"
$ws.Range("B25").Value = 495226
$ws.Range("C25").Value = "Q3"
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 24
$ws.Range("F25").Value = 27
$ws.Range("A26").Value = "495227 Q3"
$ws.Range("B26").Value = 495227
$ws.Range("C26").Value = "Q3"
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 26
$ws.Range("A27").Value = "495227 Q3"
$ws.Range("B27").Value = 495227
$ws.Range("C27").Value = "Q3"
$ws.Range("D27").Value = 24
$ws.Range("E27").Value = 26
$ws.Range("F27").Value = 26
